{"js": "// Update the statistical results table: the analysis was rerun after\n// excluding a participant who took part a second time under a changed\n// name. Replace each affected cell's old value with its recomputed\n// value. (The ADHDvCOMP column values were unaffected and are left\n// untouched.)\n\nconst replacements = [\n  [\"3.564*\", \"3.180*\"],\n  [\"-0.875\", \"-0.835\"],\n  [\"9.026*\", \"8.253*\"],\n  [\"2.723\", \"2.887\"],\n  [\"17.909*\", \"17.476*\"],\n  [\"7.485*\", \"6.906*\"],\n  [\"6.921*\", \"6.432*\"],\n  [\"-1.098\", \"-1.116\"],\n  [\"22.831*\", \"21.815*\"],\n  [\"23.558*\", \"22.505*\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchWildcards: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the statistical results table: the analysis was rerun after\n# excluding a participant who took part a second time under a changed\n# name. Replace each affected cell's old value with its recomputed value.\n# (The ADHDvCOMP column values were unaffected and are left untouched.)\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"3.564*\";  New = \"3.180*\"  },\n    @{ Old = \"-0.875\";  New = \"-0.835\"  },\n    @{ Old = \"9.026*\";  New = \"8.253*\"  },\n    @{ Old = \"2.723\";   New = \"2.887\"   },\n    @{ Old = \"17.909*\"; New = \"17.476*\" },\n    @{ Old = \"7.485*\";  New = \"6.906*\"  },\n    @{ Old = \"6.921*\";  New = \"6.432*\"  },\n    @{ Old = \"-1.098\";  New = \"-1.116\"  },\n    @{ Old = \"22.831*\"; New = \"21.815*\" },\n    @{ Old = \"23.558*\"; New = \"22.505*\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.MatchWildcards = $false\n    $found = $find.Execute()\n    if ($found) {\n        $find.Parent.Text = $pair.New\n    } else {\n        Write-Output \"NOT FOUND: $($pair.Old)\"\n    }\n}\n"}
